# "Cleaning up files and moving to .tex include"
#
# The author's last text edit in this revision landed right after
# "A paragraph of text with s|ome " (splitting what used to be one run
# of text into two). Word always tracks the location of the most
# recent edit with a hidden bookmark named "_GoBack"; since a document
# can only ever have one "_GoBack" bookmark, (re-)adding it here both
# drops the stale one (previously sitting after the second "Header 1")
# and plants the new one at the edit point - splitting the run exactly
# the way Word's own save logic does, and renumbering every other
# bookmark's w:id along the way.
$d = $word.ActiveDocument

$editPoint = $d.Content
$editPoint.Find.Execute("A paragraph of text with s", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$editPoint.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $editPoint)

# The "Normal" style now reserves 16pt (320 twips) of space after each
# paragraph.
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 16
